# Adicionado UF à Tabela de Classificação de Estados
# Rename the first column header of the "Estado" table from
# "EstadoFromCSV" to "UF" (cell A1 on the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "UF"

# Restore the active selection to the default cell (A1), since the
# original file had an explicit selection on A9 that is no longer
# present after the edit.
$ws.Range("A1").Select()
